$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1322
$ws.Range("G2").Value = 60
$ws.Range("F3").Value = 78
$ws.Range("F4").Value = 76
$ws.Range("F5").Value = 167
$ws.Range("F6").Value = 389
$ws.Range("F7").Value = 175
$ws.Range("F8").Value = 122
$ws.Range("F9").Value = 997
$ws.Range("F10").Value = 323
$ws.Range("F12").Value = 41
$ws.Range("F14").Value = 359
$ws.Range("F15").Value = 344
$ws.Range("F16").Value = 752
$ws.Range("F17").Value = 133
$ws.Range("F19").Value = 252
$ws.Range("F20").Value = 64
$ws.Range("F21").Value = 970
$ws.Range("F22").Value = 430
$ws.Range("F23").Value = 241
$ws.Range("F24").Value = 77
$ws.Range("F25").Value = 353
$ws.Range("F26").Value = 24

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 353
$ws.Range("F10").Value = 626

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1322
$ws.Range("G3").Value = 60
$ws.Range("F4").Value = 78
$ws.Range("F6").Value = 76
$ws.Range("F7").Value = 167
$ws.Range("F8").Value = 389
$ws.Range("F9").Value = 175
$ws.Range("F10").Value = 122
$ws.Range("F11").Value = 997
$ws.Range("F12").Value = 323
$ws.Range("F15").Value = 41
$ws.Range("F16").Value = 353
$ws.Range("F19").Value = 359
$ws.Range("F22").Value = 344
$ws.Range("F23").Value = 752
$ws.Range("F24").Value = 133
$ws.Range("F26").Value = 252
$ws.Range("F27").Value = 64
$ws.Range("F28").Value = 970
$ws.Range("F29").Value = 430
$ws.Range("F32").Value = 241
$ws.Range("F33").Value = 77
$ws.Range("F34").Value = 353
$ws.Range("F35").Value = 626
$ws.Range("F37").Value = 24
